# [Kadastro App] Kayıt silindi: 11364533
# Delete the record whose "Kayıt No" (column A) is 11364533 from both the
# master "Kayitlar" sheet and its per-district mirror "Merkez İlçe" sheet.
# Deleting the entire row shifts every following row up by one, which is
# exactly what the target workbook shows (row 1321 disappears from
# "Kayitlar", row 782 disappears from "Merkez İlçe", and both sheets end
# up one row shorter).

$wb = $excel.ActiveWorkbook
$recordId = "11364533"

function Remove-RecordRow($ws, $recordId) {
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
    for ($r = 2; $r -le $lastRow; $r++) {
        $val = $ws.Cells.Item($r, 1).Text
        if ($val -eq $recordId) {
            $ws.Rows.Item($r).Delete()
            return $true
        }
    }
    return $false
}

$ws1 = $wb.Worksheets.Item("Kayitlar")
Remove-RecordRow $ws1 $recordId | Out-Null

$ws5 = $wb.Worksheets.Item("Merkez İlçe")
Remove-RecordRow $ws5 $recordId | Out-Null
